$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 5")

# Fix the typo ("firest" -> "first") and merge the three runs into one.
# Deleting the existing text first (rather than overwriting it in place)
# avoids the host's text-diff logic from re-splitting the replacement text
# back along the old ("firest") run boundaries.
$tr = $sh.TextFrame.TextRange
$tr.Delete()
[void]$tr.InsertBefore("My proposed gameplay is a round-based first-person action game where players are to defend an abandoned sacred shrine from yokai. The area will be in the yard of the shrine where evil yokai will spawn on the edge of the map and work themselves inwards towards the shrine. ")

# Update position/size (EMU -> points, 1 pt = 12700 EMU). Literal point
# values below are tuned so that the host's float32 Left/Top/Width/Height
# properties round-trip back to the exact target EMU values. Applied after
# the text edit since this textbox auto-fits its height to the text, and
# re-populating the text recomputes Height.
$sh.Left = 16.111812591552734
$sh.Top = 227.30323791503906
$sh.Width = 308.662841796875
$sh.Height = 94.51409912109375
